$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "28.079.69"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -0.05%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.816.30"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +0.79%  "

$ws.Range("E4").Value = "  +0.49%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "310.50"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.13%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +0.43%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4974"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -3.37%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3914"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -1.60%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.09855"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +25.30%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "1.107"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +0.75%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "40.95"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -0.26%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "6.429"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +2.06%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "20.56"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +0.83%  "

$ws.Range("E14").Value = "  +0.49%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "1.812.14"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +1.57%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "7.267"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -0.17%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.00001138"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +5.33%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "92.21"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +0.14%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.06648"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +1.85%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +0.32%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "17.18"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -0.16%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "5.957"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +0.05%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "28.139.57"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -0.06%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "11.23"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +1.52%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.243"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +0.88%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "159.38"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -0.52%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "20.75"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +1.26%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "2.020.56"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +1.02%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "2.402"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +0.73%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "126.53"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -0.41%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "0.1056"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -2.21%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "1.033"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -1.34%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "5.570"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +0.69%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "3.622"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -0.32%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.06664"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -6.31%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.02343"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +1.39%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "8.921"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -1.17%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.2139"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +0.01%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "4.961"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -1.39%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "11.34"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -2.34%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.6199"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +0.70%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "1.189"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +2.69%  "

$ws.Range("E43").Value = "  +0.24%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "13.23"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -0.26%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.5898"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -1.24%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "3.695"
$cell.Style = "Normal"

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "1.276"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -3.04%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "124.41"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -0.99%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "1.940"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +0.86%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "1.180"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -2.56%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.06781"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -0.98%  "
